# Update "想去人数" (want-to-go count) values on the 展览 (F3,F4,F6)
# and 全部类型 (F5,F6,F8) sheets to reflect the refreshed data pull.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 356
$wsExhibit.Range("F4").Value = 2949
$wsExhibit.Range("F6").Value = 617

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 356
$wsAll.Range("F6").Value = 2949
$wsAll.Range("F8").Value = 617
